$wb = $excel.ActiveWorkbook

# ---- Metrics sheet: update raw data values ----
$wsMetrics = $wb.Worksheets.Item("Metrics")

$wsMetrics.Range("B2").Value = 397691.4
$wsMetrics.Range("B3").Value = 300975.39000000007
$wsMetrics.Range("B4").Value = 104644.44
$wsMetrics.Range("B5").Value = 16231
$wsMetrics.Range("B6").Value = 397691.4
$wsMetrics.Range("B7").Value = 300975.39000000007
$wsMetrics.Range("B8").Value = 104644.44
$wsMetrics.Range("B9").Value = 16231
$wsMetrics.Range("B10").Value = 34498943.119999997
$wsMetrics.Range("B11").Value = 32346968.18
$wsMetrics.Range("B12").Value = 12050458.299999999
$wsMetrics.Range("B13").Value = 1334138

# Recalculate the whole workbook so formulas on other sheets
# (e.g. "today") pick up the new Metrics values.
$excel.Calculate()

# Update the selection on the Metrics sheet.
$wsMetrics.Activate()
$wsMetrics.Range("D9").Select()

# ---- today sheet: selection only; A1's "=TODAY()-1" cached value ----
# updates on its own via the normal recalculation that happens after this
# script runs (the formula itself is left untouched).
$wsToday = $wb.Worksheets.Item("today")

$wsToday.Activate()
$wsToday.Range("G19").Select()
